$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.217.13'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.661.58'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.19'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5224'
$ws.Range('E6').Value = '  -1.71%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2645'
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06286'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.81'
$ws.Range('E10').Value = '  -4.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07721'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').Value = '1.675.26'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.426'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').Value = '1.888.94'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5443'
$ws.Range('E15').Value = '  -2.31%  '
$ws.Range('D16').Value = '0.0₅8160'
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.57'
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = '26.244.21'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.633'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.90'
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.05'
$ws.Range('E22').Value = '  -2.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.047'
$ws.Range('E23').Value = '  -4.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.66'
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1225'
$ws.Range('E26').Value = '  -4.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.168'
$ws.Range('E27').Value = '  -3.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.07'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.407'
$ws.Range('E29').Value = '  -2.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06025'
$ws.Range('E30').Value = '  -4.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.277'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.580'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.255'
$ws.Range('E33').Value = '  -5.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.618'
$ws.Range('E34').Value = '  -3.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9641'
$ws.Range('E35').Value = '  -4.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.427'
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5659'
$ws.Range('E38').Value = '  -8.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01594'
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.986'
$ws.Range('E40').Value = '  -3.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8549'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').Value = '1.012.75'
$ws.Range('E43').Value = '  -7.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.24'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').Value = '1.804.07'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  +6.69%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.91'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.997'
$ws.Range('E49').Value = '  -2.10%  '
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05177'
$ws.Range('E51').Value = '  -0.67%  '
